$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure the Price/Volume columns keep their original text formatting
# (avoids Excel auto-converting numeric-looking strings to real numbers).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '28.499.17'
$ws.Range("E2").Value = '  +0.42%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.869.40'
$ws.Range("E3").Value = '  +0.11%  '
$ws.Range("E4").Value = '  -1.28%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.25'
$ws.Range("E5").Value = '  -0.87%  '
$ws.Range("E6").Value = '  -1.27%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5065'
$ws.Range("E7").Value = '  -0.76%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3902'
$ws.Range("E8").Value = '  -1.60%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.08332'
$ws.Range("E9").Value = '  -0.32%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '42.83'
$ws.Range("E10").Value = '  +2.17%  '
$ws.Range("E11").Value = '  -0.40%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '6.183'
$ws.Range("E12").Value = '  -0.81%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '1.871.11'
$ws.Range("E13").Value = '  +2.51%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.27'
$ws.Range("E14").Value = '  -0.79%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '7.241'
$ws.Range("E15").Value = '  +0.53%  '
$ws.Range("E16").Value = '  -1.27%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '93.09'
$ws.Range("E17").Value = '  +2.77%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.00001097'
$ws.Range("E18").Value = '  -0.85%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06714'
$ws.Range("E19").Value = '  -0.53%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '17.60'
$ws.Range("E20").Value = '  -0.40%  '
$ws.Range("E21").Value = '  -1.26%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.909'
$ws.Range("E22").Value = '  -0.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '28.532.10'
$ws.Range("E23").Value = '  +0.56%  '
$ws.Range("E24").Value = '  -0.76%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.188'
$ws.Range("E25").Value = '  -4.36%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.078.84'
$ws.Range("E26").Value = '  +1.99%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '157.99'
$ws.Range("E27").Value = '  -2.39%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '20.54'
$ws.Range("E28").Value = '  -0.99%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.411'
$ws.Range("E29").Value = '  +2.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '125.92'
$ws.Range("E30").Value = '  -0.98%  '
$ws.Range("E31").Value = '  -1.26%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.039'
$ws.Range("E32").Value = '  +0.58%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.770'
$ws.Range("E33").Value = '  -0.20%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '3.626'
$ws.Range("E34").Value = '  -0.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.02438'
$ws.Range("E35").Value = '  +0.80%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.06551'
$ws.Range("E36").Value = '  +1.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '8.982'
$ws.Range("E37").Value = '  +1.70%  '
$ws.Range("E38").Value = '  -1.26%  '
$ws.Range("E39").Value = '  +0.60%  '
$ws.Range("E40").Value = '  +0.80%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.237'
$ws.Range("E41").Value = '  -2.18%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.6353'
$ws.Range("E42").Value = '  -0.47%  '
$ws.Range("E43").Value = '  -1.40%  '
$ws.Range("E44").Value = '  -1.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.5972'
$ws.Range("E45").Value = '  -0.82%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '13.05'
$ws.Range("E46").Value = '  +0.48%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.674'
$ws.Range("E47").Value = '  -1.01%  '
$ws.Range("E48").Value = '  +0.68%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.208'
$ws.Range("E49").Value = '  +0.57%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '122.12'
$ws.Range("E50").Value = '  +0.48%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.177'
$ws.Range("E51").Value = '  -3.03%  '
